$d = $word.ActiveDocument
$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Change 1 -------------------------------------------------------
# "Available on the start: Syringe (with some sleeping liquid inside), this
# syringe can be " paragraph gets reworked into "Available from the start: "
# (split across 3 styled runs) + a bold space + the Syringe sentence split
# into two runs ("...can be" / " placed stuck unto the tree") + a
# relocated _GoBack bookmark.
$p1 = $d.Paragraphs.Item(70)
$r1 = $p1.Range
$p1Start = $r1.Start

$xml1 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="474F21AC" w14:textId="77777777" w:rsidR="00343B16" w:rsidRPr="00D9109D" w:rsidRDefault="00343B16" w:rsidP="00343B16">' +
    '<w:r><w:t xml:space="preserve">Available </w:t></w:r>' +
    '<w:r><w:t>from</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the start:</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Syringe (with some sleeping liquid inside), this syringe can be</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> placed stuck unto the tree</w:t></w:r>' +
    '<w:bookmarkStart w:id="50" w:name="_GoBack"/><w:bookmarkEnd w:id="50"/>' +
    '</w:p>'
$r1.InsertXML($xml1) | Out-Null

# Apply the "heading 2" linked character style to the first three runs
# ("Available " + "from" + " the start:" == 25 characters).
$styleRange1 = $d.Range($p1Start, $p1Start + 25)
$styleRange1.Style = "Ttulo2Car"

# --- Change 2 -------------------------------------------------------
# The "We then see..." paragraph is split (within the same paragraph) into
# two runs, moving the page-break marker from the start of the next
# paragraph onto the second run here; the next paragraph loses that marker.
$p2 = $d.Paragraphs.Item(77)
$r2 = $p2.Range
$apostrophe = [char]0x2019
$xml2 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="29B2FB3F" w14:textId="77777777" w:rsidR="00343B16" w:rsidRDefault="00343B16" w:rsidP="00343B16">' +
    '<w:r><w:t xml:space="preserve">We then see what they meant when they said this. In the middle of the room are 3 recipients (jugs?) of different sizes on a table. The jugs carry 8, 5, and 3 Liters respectively, and the player is instructed to </w:t></w:r>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>measure 4 liters to dilute the sleeping agent from the syringe, as it is otherwise lethal. This puzzle in specific mirrors Taylor' + $apostrophe + 's actions.</w:t></w:r>' +
    '</w:p>'
$r2.InsertXML($xml2) | Out-Null

$p3 = $d.Paragraphs.Item(78)
$r3 = $p3.Range
$xml3 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="525D5E00" w14:textId="77777777" w:rsidR="00343B16" w:rsidRDefault="00343B16" w:rsidP="00343B16">' +
    '<w:r><w:t>If you consider the jugs will be placed in the order 8L, 5L, and 3L, with only the 8L jug being full of wine the solution goes as follows:</w:t></w:r>' +
    '</w:p>'
$r3.InsertXML($xml3) | Out-Null

# --- Change 3 -------------------------------------------------------
# The trailing paragraph that used to hold the old _GoBack bookmark becomes
# a plain empty paragraph.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$xmlLast = '<w:p ' + $wNs + '/>'
$rLast.InsertXML($xmlLast) | Out-Null
